$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 16, shifting the existing row 16 (and below) down to row 17.
$ws.Rows.Item(16).Insert()

# New row 16: weekly price entry for date 45021 (2023-04-05)
$ws.Cells.Item(16, 1).Value = 1
$ws.Cells.Item(16, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(16, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(16, 4).Value = 45021
$ws.Cells.Item(16, 4).NumberFormat = $ws.Cells.Item(17, 4).NumberFormat
$ws.Cells.Item(16, 5).Value = 15
$ws.Cells.Item(16, 6).Value = "Fruta"
$ws.Cells.Item(16, 7).Value = 100107
$ws.Cells.Item(16, 8).Value = "Otros"
$ws.Cells.Item(16, 9).Value = 100107011
$ws.Cells.Item(16, 10).Value = "Tuna"
$ws.Cells.Item(16, 11).Value = "Sin especificar"
$ws.Cells.Item(16, 12).Value = "Segunda"
$ws.Cells.Item(16, 13).Value = 250
$ws.Cells.Item(16, 14).Value = 22000
$ws.Cells.Item(16, 15).Value = 23000
$ws.Cells.Item(16, 16).Value = 22500
$ws.Cells.Item(16, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item(16, 18).Value = "Región de Coquimbo"
$ws.Cells.Item(16, 19).Value = 1125
$ws.Cells.Item(16, 20).Value = 20
